$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '54.411.24'
$ws.Range("E2").Value = '  -2.47%  '
$ws.Range("D3").Value = '2.286.54'
$ws.Range("E3").Value = '  -2.89%  '
$ws.Range("E4").Value = '  -0.26%  '
$ws.Range("D5").Value = "'493.56"
$ws.Range("E5").Value = '  -1.86%  '
$ws.Range("D6").Value = "'127.20"
$ws.Range("E6").Value = '  -2.12%  '
$ws.Range("D7").Value = "'0.999"
$ws.Range("E7").Value = '  +0.10%  '
$ws.Range("D8").Value = "'0.529"
$ws.Range("E8").Value = '  -1.66%  '
$ws.Range("D9").Value = '2.286.57'
$ws.Range("E9").Value = '  -3.45%  '
$ws.Range("D10").Value = "'0.0947"
$ws.Range("E10").Value = '  -1.84%  '
$ws.Range("E11").Value = '  +0.18%  '
$ws.Range("D12").Value = "'0.323"
$ws.Range("E12").Value = '  +0.89%  '
$ws.Range("E13").Value = '  -3.34%  '
$ws.Range("D14").Value = '2.673.82'
$ws.Range("E14").Value = '  -3.69%  '
$ws.Range("D15").Value = "'21.60"
$ws.Range("E15").Value = '  +0.66%  '
$ws.Range("D16").Value = '54.320.23'
$ws.Range("E16").Value = '  -2.57%  '
$ws.Range("E17").Value = '  -2.19%  '
$ws.Range("D18").Value = '2.271.72'
$ws.Range("E18").Value = '  -4.80%  '
$ws.Range("D19").Value = "'10.00"
$ws.Range("E19").Value = '  +1.19%  '
$ws.Range("D20").Value = "'4.06"
$ws.Range("E20").Value = '  +1.41%  '
$ws.Range("D21").Value = "'303.46"
$ws.Range("E21").Value = '  -2.06%  '
$ws.Range("D22").Value = "'6.48"
$ws.Range("E22").Value = '  +4.38%  '
$ws.Range("D23").Value = "'0.999"
$ws.Range("E23").Value = '  +0.08%  '
$ws.Range("E24").Value = '  -3.47%  '
$ws.Range("D25").Value = "'63.44"
$ws.Range("E25").Value = '  -2.78%  '
$ws.Range("D26").Value = "'1.01"
$ws.Range("E26").Value = '  +0.85%  '
$ws.Range("E27").Value = '  +1.08%  '
$ws.Range("D28").Value = '2.376.43'
$ws.Range("E28").Value = '  -3.81%  '
$ws.Range("D29").Value = "'0.150"
$ws.Range("E29").Value = '  +3.60%  '
$ws.Range("D30").Value = "'7.09"
$ws.Range("E30").Value = '  -0.53%  '
$ws.Range("D31").Value = "'168.65"
$ws.Range("E31").Value = '  -1.19%  '
$ws.Range("D32").Value = "'1.59"
$ws.Range("E32").Value = '  -2.19%  '
$ws.Range("D33").Value = '0.0₃0685'
$ws.Range("E33").Value = '  -2.83%  '
$ws.Range("E34").Value = '  +2.50%  '
$ws.Range("E35").Value = '  -0.13%  '
$ws.Range("E36").Value = '  +0.28%  '
$ws.Range("D37").Value = "'1.08"
$ws.Range("E37").Value = '  +1.18%  '
$ws.Range("D38").Value = "'17.62"
$ws.Range("E38").Value = '  -0.06%  '
$ws.Range("E39").Value = '  +2.54%  '
$ws.Range("D40").Value = "'0.866"
$ws.Range("E40").Value = '  +3.11%  '
$ws.Range("D41").Value = "'3.64"
$ws.Range("E41").Value = '  -0.01%  '
$ws.Range("D42").Value = "'35.50"
$ws.Range("E42").Value = '  -1.82%  '
$ws.Range("D43").Value = "'0.374"
$ws.Range("E43").Value = '  +0.92%  '
$ws.Range("D44").Value = "'1.40"
$ws.Range("E44").Value = '  +0.47%  '
$ws.Range("D45").Value = "'128.94"
$ws.Range("E45").Value = '  +2.38%  '
$ws.Range("E46").Value = '  +0.14%  '
$ws.Range("D47").Value = "'4.87"
$ws.Range("E47").Value = '  +0.16%  '
$ws.Range("E48").Value = '  -0.44%  '
$ws.Range("D49").Value = "'0.547"
$ws.Range("E49").Value = '  -1.70%  '
$ws.Range("D50").Value = "'238.79"
$ws.Range("E50").Value = '  -0.33%  '
$ws.Range("D51").Value = "'0.0480"
$ws.Range("E51").Value = '  +0.66%  '
